$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts old B,C,D -> C,D,E)
$ws.Range("B1").EntireColumn.Insert()

# Fill in the new column B values
$ws.Range("B2").Value = 3
$ws.Range("B3").Value = "id2"
$ws.Range("B4").Value = "INT"
$ws.Range("B5").Value = 11
$ws.Range("B6").Value = 12
$ws.Range("B7").Value = 13
$ws.Range("B8").Value = 14

# Column A4 text changed from INT to key
$ws.Range("A4").Value = "key"

# D column (was C) values updated
$ws.Range("D7").Value = 1600
$ws.Range("D8").Value = 2200

# New block starting row 11 - KR region header
$ws.Range("A11:B11").Interior.Color = 255
$ws.Range("A11").Value = "KR"

$ws.Range("A6:E6").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Range("A12").Value = 1
$ws.Range("B12").Value = 12
$ws.Range("C12").Value = 6
$ws.Range("D12").Value = 301
$ws.Range("E12").Value = "[{1,2}]"

# New block starting row 15 - TW region header
$ws.Range("A15:B15").Interior.Color = 255
$ws.Range("A15").Value = "TW"

$ws.Range("A5:E5").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Range("A16").Value = 0
$ws.Range("B16").Value = 11
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 62
$ws.Range("E16").Value = "[1,2]"

# New block starting row 19 - Th region header
$ws.Range("A19").Interior.Color = 255
$ws.Range("A19").Value = "Th"

# Restore the selection to match the target workbook view
$ws.Range("B5").Select()
